# Adds June grocery data (rows 102-149) plus a few minor tweaks,
# matching commit "Added June data plus a few minor tweaks."

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1) Pre-seed the shared-strings table so brand-new item names land in
#    the exact order they appear in the target workbook. Excel allocates
#    shared-string ids in first-use order, and the true edit history here
#    didn't touch every new cell strictly top-to-bottom, so we stage the
#    values (in the desired final order) in a scratch cell far off the
#    used range, then clear them before writing real data.
# ---------------------------------------------------------------------
$newItemOrder = @(
    'Peeled tomatoes',
    'Ginger',
    'Toothpaste',
    'Flowers',
    'Kiwi fruit',
    'Chicken',
    'Spring onions',
    'Yoghurt',
    'Coffee',
    'Only organic breakfast',
    'Baby wipes',
    'Sweet potato',
    'Vegemite',
    'Hand sanitiser',
    'Only organic beef pasta',
    'Only organic cheese pasta',
    'Pumpkin',
    'Bananas',
    'Only organic rice cakes',
    'Courgette',
    'Mandarins',
    'toffee',
    'Kiwi garden yoghurt drops',
    'Band-Aid adhesive strips',
    'Licorice'
)

$scratch = $ws.Range("ZZ1")
foreach ($name in $newItemOrder) {
    $scratch.Value2 = $name
    $scratch.ClearContents()
}

# ---------------------------------------------------------------------
# 2) New transaction rows (June data), in final row order 102-149.
#    Columns: item, amount, date (serial).
# ---------------------------------------------------------------------
$newRows = @(
    @('Peeled tomatoes', 5.6, 43646),
    @('Tomato paste', 4.5, 43646),
    @('Ginger', 0.9, 43646),
    @('Brown onions', 1.7, 43646),
    @('Flowers', 9.99, 43646),
    @('Toothpaste', 2.5, 43646),
    @('Chocolate', 2.9, 43645),
    @('Chocolate', 1.7, 43645),
    @('Tomatoes', 1.4, 43645),
    @('Kiwi fruit', 1.3, 43645),
    @('Chicken', 7.5, 43645),
    @('Spring onions', 2.8, 43645),
    @('Milk', 3.4, 43643),
    @('Milk', 3.4, 43643),
    @('Yoghurt', 8.4, 43643),
    @('Coffee', 12, 43643),
    @('Only organic breakfast', 4.5999999999999996, 43643),
    @('Bread', 8.4, 43643),
    @('Chocolate', 5.3, 43642),
    @('Tomato paste', 3, 43641),
    @('Baby wipes', 6, 43641),
    @('Sweet potato', 13, 43641),
    @('Vegemite', 8.5, 43641),
    @('Hand sanitiser', 10, 43641),
    @('Only organic beef pasta', 3.5, 43641),
    @('Only organic cheese pasta', 3.5, 43641),
    @('Only organic breakfast', 4.5999999999999996, 43641),
    @('Only organic breakfast', 2.2999999999999998, 43641),
    @('Pumpkin', 3, 43641),
    @('Brown onions', 0.3, 43641),
    @('Bananas', 3.5, 43641),
    @('Only organic rice cakes', 4.3, 43641),
    @('Courgette', 1.4, 43641),
    @('Mandarins', 1.7, 43641),
    @('Free range eggs', 7.6, 43641),
    @('Brown onions', 1.6, 43639),
    @('Chocolate', 3, 43639),
    @('toffee', 2, 43639),
    @('Milk', 3.4, 43632),
    @('Milk', 3.4, 43632),
    @('Kiwi garden yoghurt drops', 4, 43634),
    @('Coffee', 6, 43634),
    @('Band-Aid adhesive strips', 5.7, 43634),
    @('Bananas', 2.25, 43634),
    @('Licorice', 1.2, 43634),
    @('Chocolate', 4.2, 43634),
    @('Toothpaste', 4, 43634),
    @('Chocolate', 1, 43634)
)

$startRow = 102
$r = $startRow
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value2 = $row[0]
    $ws.Cells.Item($r, 2).Value2 = $row[1]
    $ws.Cells.Item($r, 3).Value2 = $row[2]
    $r = $r + 1
}
$endRow = $r - 1

# ---------------------------------------------------------------------
# 3) Formatting: reuse the existing styles used throughout the sheet
#    (column A = wrap-off item-name style, column C = date style) by
#    copying formats from the last pre-existing data row (101).
# ---------------------------------------------------------------------
$ws.Range("A101").Copy() | Out-Null
$ws.Range("A102:A$endRow").PasteSpecial(-4122) | Out-Null

$ws.Range("C101").Copy() | Out-Null
$ws.Range("C102:C$endRow").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

# Row 149's amount cell carries the date-style's sibling formatting
# (same font style as column A/C) in the source workbook.
$ws.Range("B3").Copy() | Out-Null
$ws.Range("B149").PasteSpecial(-4122) | Out-Null
$ws.Range("B149").Value2 = 1

# Row 137 has a stray formatted (but empty) D column cell.
$ws.Range("C101").Copy() | Out-Null
$ws.Range("D137").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 4) Column D width (new column introduced alongside the stray D137
#    cell). Target best-fit width is ~9.43 chars; nudge ColumnWidth to
#    land on the closest width bucket this engine will serialize.
# ---------------------------------------------------------------------
$ws.Columns.Item(4).ColumnWidth = 8.6

# ---------------------------------------------------------------------
# 5) Restore the view state (selection / scroll position) to reflect
#    the newly-added bottom of the data range.
# ---------------------------------------------------------------------
$ws.Range("C150").Select()
$excel.ActiveWindow.ScrollRow = 108
$excel.ActiveWindow.ScrollColumn = 1
